# Employee.xlsx update: replace the 5-employee sample data with a new
# 4-employee sample data set (john doe, tracy, oliver queen, bhavesh),
# and remove the now-unused 5th data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - John Doe
$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "john doe"
$ws.Range("C2").Value = "10 clover lane"
$ws.Range("D2").Value = "john@gmail.com"
$ws.Range("E2").Value = "01/09/2020"
$ws.Range("F2").Value = 784521
$ws.Range("G2").Value = "ux designer"

# Row 3 - Tracy
$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "tracy"
$ws.Range("C3").Value = "10 field lane"
$ws.Range("D3").Value = "tracy@gmail.com"
$ws.Range("E3").Value = "11/08/2020"
$ws.Range("F3").Value = 54875421
$ws.Range("G3").Value = "java developer"

# Row 4 - Oliver Queen
$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "oliver queen"
$ws.Range("C4").Value = "starline city"
$ws.Range("D4").Value = "oliver@gmail.com"
$ws.Range("E4").Value = "10/30/2020"
$ws.Range("F4").Value = 78545212
$ws.Range("G4").Value = "a full stack developer"

# Row 5 - Bhavesh
$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "bhavesh"
$ws.Range("C5").Value = "10 front lane"
$ws.Range("D5").Value = "bhavesh@gmai.com"
$ws.Range("E5").Value = "11/05/2020"
$ws.Range("F5").Value = 784512
$ws.Range("G5").Value = "ux designer"

# Row 6 is no longer part of the sample data set - remove it entirely.
$ws.Rows.Item(6).Delete()
